$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.150782346725464
$ws.Range("B1").Value = 2.253428936004639
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.145925998687744
$ws.Range("E1").Value = 1.063439130783081
